$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "24/10/2025"
$ws.Range("B13").Value = "Termalica B-B."
$ws.Range("C13").Value = 1
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = "Zaglebie"
$ws.Range("F13").Value = "D"
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 1.08
$ws.Range("L13").Value = 1.43
$ws.Range("M13").Value = 13
$ws.Range("N13").Value = 11
$ws.Range("O13").Value = 6
$ws.Range("P13").Value = 3
